$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "GarProVen",
    "GrexActFij",
    "GrexMsgs",
    "grexrepo",
    "HisManAct",
    "HVActivo",
    "InvActDep",
    "InvActPer.sql",
    "InvCom1",
    "InvComp",
    "InveDepTpl",
    "LisConDepXlsDet",
    "LisConDepXls",
    "LisDepMesXlsDet",
    "LisDepMesXlsRes",
    "LISTEREJE",
    "LisTraslPer",
    "LisTrasl",
    "logger_asserts",
    "NumerosErroresActivosFijos",
    "Persempr",
    "PolManAct111",
    "POLMANACT",
    "PolManAct",
    "PolSegAct",
    "reportes_acfi",
    "SinPolMan",
    "subActAdi2",
    "subActAdi",
    "subActBaj",
    "subActDes",
    "subActImg",
    "subActPer",
    "subActTras",
    "TrasActPer"
)

# Clear the old sample rows (rows 2-4) before writing the real data.
$ws.Range("A2:B4").ClearContents()

$row = 2
foreach ($name in $names) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = 29
    $row = $row + 1
}

$ws.Columns.Item(1).ColumnWidth = 52.28515625

$ws.Range("A36").Select()
$excel.ActiveWindow.ScrollRow = 22
